$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A width and number format (date/time) for the whole column
$colA = $ws.Range("A1:A3")
$colA.NumberFormat = "m/d/yy h:mm"

# Row 2 data
$ws.Range("A2").Value = Get-Date -Year 2016 -Month 8 -Day 27 -Hour 15 -Minute 25 -Second 10
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 55
$ws.Range("D2").Value = 43
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 14426
$ws.Range("H2").Value = 10760
$ws.Range("I2").Value = 576
$ws.Range("J2").Value = 96
$ws.Range("K2").Value = 75
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "Named"

# Row 3 data
$ws.Range("A3").Value = Get-Date -Year 2016 -Month 8 -Day 27 -Hour 15 -Minute 31 -Second 48
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 55
$ws.Range("D3").Value = 43
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 13948
$ws.Range("H3").Value = 10709
$ws.Range("I3").Value = 568
$ws.Range("J3").Value = 95
$ws.Range("K3").Value = 75
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "Named"

# Resize column A to fit the new date/time values (matches Excel's recalculated best-fit width)
$ws.Columns.Item(1).ColumnWidth = 14
